$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the "Image Path" column (J) for rows 47-115 ---
# Each row gets an image filename derived from its core name (column A),
# e.g. "ELP 22/6/16" -> "images/elp22616.png". Entered in the same order
# the original author typed them (mostly top-to-bottom, with a couple of
# out-of-order corrections), so the shared-string table is built up
# identically.
$imagePaths = @(
    @(49, "images/elp22616.png"),
    @(50, "images/elp32620.png"),
    @(51, "images/elp38825.png"),
    @(52, "images/elp431028.png"),
    @(53, "images/elp431028.png"),
    @(54, "images/elp641050.png"),
    @(55, "images/elp1022038.png"),
    @(56, "images/er955.png"),
    @(57, "images/er115.png"),
    @(58, "images/er1456.png"),
    @(59, "images/er18310.png"),
    @(60, "images/er23513.png"),
    @(61, "images/er23513.png"),
    @(62, "images/er25615.png"),
    @(63, "images/er25615.png"),
    @(64, "images/er281711.png"),
    @(65, "images/er32521.png"),
    @(48, "images/elp18410.png"),
    @(47, "images/elp14355.png"),
    @(66, "images/er352011.png"),
    @(67, "images/er422215.png"),
    @(68, "images/er461718.png"),
    @(69, "images/er492717.png"),
    @(70, "images/er541818.png"),
    @(71, "images/eq1328587.png"),
    @(72, "images/i13187.png"),
    @(73, "images/eq206314.png"),
    @(74, "images/i202314.png"),
    @(76, "images/i252318.png"),
    @(77, "images/eq30820.png"),
    @(78, "images/i302720.png"),
    @(75, "images/eq255618.png"),
    @(79, "images/rm4.png"),
    @(80, "images/rm5.png"),
    @(81, "images/rm6.png"),
    @(82, "images/rm7.png"),
    @(83, "images/rm8.png"),
    @(84, "images/rm10.png"),
    @(85, "images/rm12.png"),
    @(86, "images/rm14.png"),
    @(87, "images/pm5039.png"),
    @(88, "images/pm6249.png"),
    @(89, "images/pm7459.png"),
    @(90, "images/pm8770.png"),
    @(91, "images/pm11493.png"),
    @(92, "images/p3326.png"),
    @(93, "images/p4641.png"),
    @(94, "images/p5833.png"),
    @(95, "images/p74.png"),
    @(96, "images/p95.png"),
    @(97, "images/p117.png"),
    @(98, "images/p148.png"),
    @(99, "images/p1811.png"),
    @(100, "images/p2213.png"),
    @(101, "images/p2616.png"),
    @(102, "images/p3019.png"),
    @(103, "images/p3622.png"),
    @(104, "images/p4125.png"),
    @(105, "images/p4728.png"),
    @(106, "images/p16116.png"),
    @(107, "images/p2016.png"),
    @(108, "images/p2020.png"),
    @(109, "images/p2620.png"),
    @(110, "images/p2625.png"),
    @(111, "images/p3220.png"),
    @(112, "images/p3230.png"),
    @(113, "images/p3535.png"),
    @(114, "images/p4040.png"),
    @(115, "images/p5050.png")

)

foreach ($entry in $imagePaths) {
    $r = $entry[0]
    $path = $entry[1]
    $ws.Cells.Item($r, 10).Value = $path
}

# --- Header row tidy-up ---
# J1 picks up the same style already used by H1 (border/left-align/wrap),
# replacing its old one-off style.
$ws.Range("H1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 1 no longer needs the taller custom height.
$ws.Rows.Item(1).RowHeight = 15

# --- View state ---
# Selection moves to K114 (scrolled down near the bottom of the table).
$ws.Range("K114").Select()
